# tambah fitur upload siswa excel
#
# Re-capitalize the student header labels, insert a new "NO" numbering
# column at the front of the sheet, and tidy the header row's
# font size / column widths / row height accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing headers in place first so their shared-string order
# stays stable: NIS, NIK, Nama_Siswa, Alamat, Gender, Kontak, Kelas, Nama_Ibu
$ws.Range("A1").Value = "NIS"
$ws.Range("B1").Value = "NIK"
$ws.Range("C1").Value = "Nama_Siswa"
$ws.Range("D1").Value = "Alamat"
$ws.Range("E1").Value = "Gender"
$ws.Range("F1").Value = "Kontak"
$ws.Range("G1").Value = "Kelas"
$ws.Range("H1").Value = "Nama_Ibu"

# Insert a brand new column A; the headers above shift right from A..H to B..I
$ws.Columns("A").Insert()

# Give the new A1 the same bold header style as the rest of the row, then
# fill in its text ("NO" becomes the newest shared string, appended last).
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "NO"

# Shrink the header font from 14pt to 12pt
$ws.Range("A1:I1").Font.Size = 12

# Widen the new NIK column (C) and tidy the header row height for the
# smaller font
$ws.Columns("C").ColumnWidth = 28.6
$ws.Rows(1).RowHeight = 15.75

# Match the active cell selection saved in the workbook
$ws.Range("K6").Select() | Out-Null
